$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44316
$ws.Range("M2").Value2 = 120

# Row 3
$ws.Range("D3").Value2 = 44313
$ws.Range("M3").Value2 = 120

# Row 4
$ws.Range("D4").Value2 = 44323
$ws.Range("M4").Value2 = 80

# Row 5
$ws.Range("D5").Value2 = 44302
$ws.Range("M5").Value2 = 80

# Row 6
$ws.Range("D6").Value2 = 44330
$ws.Range("M6").Value2 = 60

# Row 7
$ws.Range("D7").Value2 = 44327

# Row 8
$ws.Range("D8").Value2 = 44306
$ws.Range("M8").Value2 = 80

# Row 9
$ws.Range("D9").Value2 = 44322
$ws.Range("M9").Value2 = 60
$ws.Range("Q9").Value2 = "$/caja 10 kilos empedrada"
$ws.Range("S9").Value2 = 11500
$ws.Range("T9").Value2 = 1

# Row 10
$ws.Range("D10").Value2 = 44309
$ws.Range("M10").Value2 = 80
$ws.Range("Q10").Value2 = "$/caja 14 kilos granel"
$ws.Range("S10").Value2 = 821
$ws.Range("T10").Value2 = 14
